# "Linking settings tab to selection tab"
#
# Renames the three config tabs to their capitalized/localized display
# names and switches the active/selected tab from "german" to "english":
#   exiobase -> Exiobase
#   german   -> Deutsch   (no longer the selected tab)
#   english  -> English   (becomes the selected tab, cursor parked at A32)

$wb = $excel.ActiveWorkbook

$wsExiobase = $wb.Worksheets.Item("exiobase")
$wsGerman   = $wb.Worksheets.Item("german")
$wsEnglish  = $wb.Worksheets.Item("english")

$wsExiobase.Name = "Exiobase"
$wsGerman.Name   = "Deutsch"
$wsEnglish.Name  = "English"

# Make the English sheet the active/selected tab and move its selection
# to A32 (it previously pointed at A34 while "Deutsch" was the active tab).
$wsEnglish.Activate() | Out-Null
$wsEnglish.Range("A32").Select() | Out-Null
